$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Use" column descriptions for each BOM row (F3:F25)
$ws.Range("F3").Value = "Used to print the dispenser"

# New "Use" column header (F2) - Times New Roman to match the Part/Quantity/Cost/Link headers
$ws.Range("F2").Value = "Use"
$ws.Range("F2").Font.Name = "Times New Roman"

$ws.Range("F4").Value = "Controls the dispenser"
$ws.Range("F5").Value = "Breaks out HDMI connection on Pi to connector panel"
$ws.Range("F6").Value = "Breaks out 5V USB power on Pi to connector panel"
$ws.Range("F7").Value = "Breaks out 12V power on Pi HAT to connector panel"
$ws.Range("F8").Value = "Converts the HDMI panel connector to Micro HDMI on Pi"
$ws.Range("F9").Value = "12V power supply for stepper motor"
$ws.Range("F10").Value = "Motor to drive the treat jogger"
$ws.Range("F11").Value = "5V power supply for Raspberry Pi"
$ws.Range("F12").Value = "Acrylic plate to cover the treats from the environment"
$ws.Range("F13").Value = "External sensors to detect a treat falling from jogger"
$ws.Range("F14").Value = "Used to hold down wires within the dispenser"
$ws.Range("F15").Value = "Used to hold data and OS on Raspberry Pi"
$ws.Range("F16").Value = "Connects the panel USB C to Raspberry Pi "
$ws.Range("F17").Value = "Double pull double throw switch for both 5V and 12V power"
$ws.Range("F18").Value = "Placed on Pi HAT into the female headers"
$ws.Range("F19").Value = "Soldered to Pi HAT"
$ws.Range("F20").Value = "Soldered to Pi HAT"
$ws.Range("F21").Value = "Soldered to Pi HAT"
$ws.Range("F22").Value = "Soldered to Pi HAT"
$ws.Range("F23").Value = "Soldered to Pi HAT"
$ws.Range("F24").Value = "Soldered to Pi HAT"
$ws.Range("F25").Value = "Soldered to Pi HAT"

# Match the author's final cursor position/selection
$ws.Range("F26").Select() | Out-Null
